$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.799.53'
$ws.Range('E2').Value = '  -1.48%  '
$ws.Range('D3').Value = '1.634.19'
$ws.Range('E3').Value = '  -1.48%  '
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.60'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5032'
$ws.Range('E6').Value = '  -2.31%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2575'
$ws.Range('E8').Value = '  -0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06427'
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.64'
$ws.Range('E10').Value = '  -1.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07703'
$ws.Range('E11').Value = '  -0.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.248'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('D13').Value = '1.636.17'
$ws.Range('E13').Value = '  -1.42%  '
$ws.Range('D14').Value = '1.859.48'
$ws.Range('E14').Value = '  -1.45%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5459'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('D16').Value = '0.0₅7946'
$ws.Range('E16').Value = '  -1.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.57'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '25.823.68'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '203.43'
$ws.Range('E20').Value = '  -3.50%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.322'
$ws.Range('E21').Value = '  -2.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.957'
$ws.Range('E22').Value = '  -1.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.968'
$ws.Range('E23').Value = '  -1.16%  '
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.920'
$ws.Range('E25').Value = '  +9.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '141.21'
$ws.Range('E26').Value = '  -2.85%  '
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.71'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.713'
$ws.Range('E29').Value = '  -3.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05032'
$ws.Range('E30').Value = '  -3.15%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.241'
$ws.Range('E31').Value = '  -0.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.272'
$ws.Range('E32').Value = '  -2.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.186'
$ws.Range('E33').Value = '  -1.29%  '
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('D36').Value = '1.178.95'
$ws.Range('E36').Value = '  +0.08%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8962'
$ws.Range('E37').Value = '  -3.69%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.603'
$ws.Range('E38').Value = '  -5.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5613'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01560'
$ws.Range('E40').Value = '  -2.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.548'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('E42').Value = '  -0.33%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.671'
$ws.Range('E43').Value = '  -0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8064'
$ws.Range('E44').Value = '  -3.86%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.73'
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('D46').Value = '1.771.94'
$ws.Range('E46').Value = '  -1.42%  '
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4512'
$ws.Range('E48').Value = '  -0.60%  '
$ws.Range('E49').Value = '  -0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '54.88'
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05040'
$ws.Range('E51').Value = '  -0.43%  '
